$d = $word.ActiveDocument

# Locate the last paragraph ("Der Vermieter kann die Wohnung ohne bestimmten
# Anlass besichtigen.") and append two new paragraphs right after it.
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)

$insertRange = $lastPara.Range
$insertRange.Collapse(0)  # wdCollapseEnd

$insertRange.InsertAfter("`rEine Gebrauchsüberlassung der Mietsache an Dritte ist ausgeschlossen.`rDas Recht des Mieters, eine Mietminderung gem. § 536 BGB zu erklären, ist ausgeschlossen.")
